$p = $ppt.ActivePresentation

# Slide 2 ("Профессии:") holds the four profession entries; each of its
# four paragraphs links to the slide that details that profession.
#   paragraph 1: Web-дизайнер    -> slide 3
#   paragraph 2: QA тестировщик  -> slide 4
#   paragraph 3: Backend разработчик -> slide 5
#   paragraph 4: Frontend разработчик -> slide 6
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

$targets = @(3, 4, 5, 6)

for ($i = 0; $i -lt $targets.Length; $i++) {
    $paraIndex = $i + 1
    $destSlide = $p.Slides.Item($targets[$i])
    $para = $tr.Paragraphs($paraIndex, 1)

    $action = $para.ActionSettings(1)
    $action.Action = 7  # ppActionHyperlink
    $action.Hyperlink.SubAddress = $destSlide.SlideID
}
